# Update "想去人数" (interested-attendee counts) on the 展览 and 全部类型 sheets
# to reflect the latest generated data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 156
$ws1.Range("F4").Value = 723
$ws1.Range("F5").Value = 63

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 156
$ws4.Range("F5").Value = 723
$ws4.Range("F6").Value = 63
